$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (Coin name / Link / Volume label columns): plain assignment keeps them as text
# since none of these values parse as numbers.
$textUpdates = @{
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E10' = '9WazirXWRX'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E11' = '10LiechtensteinCryptoassetsExchangeLCX'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E13' = '12BitrueCoinBTR'
    'B14' = 'MCDex'
    'C14' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'E14' = '13MCDexMCB'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E15' = '14BitMartTokenBMX'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E16' = '15BitForexTokenBF'
    'B17' = 'CoinExToken'
    'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E17' = '16CoinExTokenCET'
    'B18' = 'One'
    'C18' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E18' = '17OneONE'
    'E27' = '26UpBotsUBXT'
    'B41' = 'BKEXToken'
    'C41' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E41' = '40BKEXTokenBKK'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'E42' = '41CEJICEJI'
    'B43' = 'KickToken'
    'C43' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'E43' = '42KickTokenKICKWorstin24h'
    'E47' = '46CoinbaseStockTokenCOINBestin24h'
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Numeric-looking cells (Price column) are stored as TEXT in the source workbook (t="inlineStr"),
# e.g. "244.54" rather than the number 244.54. Assigning a plain numeric string via .Value would
# have Excel auto-convert it to a real number, so force a Text number format before the write and
# clear the temporary formatting afterwards to avoid leaving a stray style behind.
$numericTextUpdates = @{
    'D2' = '244.54'
    'D3' = '22.02'
    'D4' = '5.388'
    'D5' = '0.05848'
    'D6' = '3.394'
    'D7' = '6.352'
    'D8' = '0.8170'
    'D9' = '1.005'
    'D10' = '0.1425'
    'D11' = '0.03617'
    'D12' = '0.07427'
    'D13' = '0.03039'
    'D14' = '4.216'
    'D15' = '0.09386'
    'D16' = '0.001605'
    'D17' = '0.04832'
    'D18' = '0.0005898'
    'D19' = '0.006047'
    'D20' = '0.004088'
    'D21' = '0.0009982'
    'D22' = '0.0001503'
    'D23' = '3.688'
    'D24' = '2.222'
    'D27' = '0.0002904'
    'D40' = '0.03856'
    'D41' = '0.1074'
    'D42' = '0.002507'
    'D43' = '0.003014'
    'D44' = '0.006244'
    'D45' = '0.00005630'
    'D47' = '0.8211'
    'D48' = '0.1424'
    'D49' = '0.00002103'
    'D50' = '0.01012'
}
foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$ref]
    $cell.ClearFormats()
}
